$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Should be at top after hello row when opening"
$ws.Range("A4").Value = "Another row"

$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
